$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1.082
$ws.Range("E2").Value = 4.098
$ws.Range("F2").Value = 1.404
$ws.Range("G2").Value = 24.492
$ws.Range("I2").Value = 0.00002965634507584613
$ws.Range("K2").Value = -45.27869493799358
$ws.Range("M2").Value = -9.623577012633623
